$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Förändrad" (changed) date in column C for all data rows (2-51)
#    from 45184 to 45186.
for ($r = 2; $r -le 51; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value = 45186
    }
}

# 2) Append the "Beteckning" (designation) text as a second HYPERLINK()
#    argument for the link columns (S, T, V, W, X, Y) on rows that have
#    those formulas (rows 2-9).
$cols = @("S", "T", "V", "W", "X", "Y")
for ($r = 2; $r -le 9; $r++) {
    $label = $ws.Cells.Item($r, 1).Value2
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$r")
        $f = $cell.Formula
        if ($f -ne $null -and $f -ne "") {
            $trimmed = $f.TrimEnd()
            if ($trimmed.EndsWith(")")) {
                $newFormula = $trimmed.Substring(0, $trimmed.Length - 1) + ', "' + $label + '")'
                $cell.Formula = $newFormula
            }
        }
    }
}
